# BRVM Realtime & Sectors
# The "Market Cap" column (K) was stored as text strings with thousands
# separators (e.g. "1,781,000,000,000"). Convert every one of those cells
# (rows 2-41) into a genuine numeric value so Excel stores/treats them as
# numbers rather than text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 1781000000000
$ws.Range("K3").Value = 1619545012500
$ws.Range("K4").Value = 550666647000
$ws.Range("K5").Value = 372967815000
$ws.Range("K6").Value = 325513920996
$ws.Range("K7").Value = 281600000000
$ws.Range("K8").Value = 280000000000
$ws.Range("K9").Value = 159307496145
$ws.Range("K10").Value = 158400000000
$ws.Range("K11").Value = 149416608000
$ws.Range("K12").Value = 145000000000
$ws.Range("K13").Value = 143209308000
$ws.Range("K14").Value = 143140000000
$ws.Range("K15").Value = 133308346500
$ws.Range("K16").Value = 131722003380
$ws.Range("K17").Value = 131711595900
$ws.Range("K18").Value = 122500024500
$ws.Range("K19").Value = 113116500000
$ws.Range("K20").Value = 104659569320
$ws.Range("K21").Value = 103571832000
$ws.Range("K22").Value = 92400000000
$ws.Range("K23").Value = 83013832500
$ws.Range("K24").Value = 80212608000
$ws.Range("K25").Value = 76680000000
$ws.Range("K26").Value = 75606440000
$ws.Range("K27").Value = 75075000000
$ws.Range("K28").Value = 74277156000
$ws.Range("K29").Value = 57761091300
$ws.Range("K30").Value = 55100400000
$ws.Range("K31").Value = 48195000000
$ws.Range("K32").Value = 39510000000
$ws.Range("K33").Value = 26992500000
$ws.Range("K34").Value = 26797106000
$ws.Range("K35").Value = 17254400000
$ws.Range("K36").Value = 14021920000
$ws.Range("K37").Value = 11101250000
$ws.Range("K38").Value = 10819200000
$ws.Range("K39").Value = 8820000000
$ws.Range("K40").Value = 7551360000
$ws.Range("K41").Value = 6008800000
